$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Through 2022-05-06 -> Through 2022-05-07)
$ws.Name = "Through 2022-05-07"

# Update the column header text to match the new "through" date
$ws.Range("B1").Value = "May 2022 (through May 07)"

# Row 2 - Englewood: new carjacking in May 2021 (column G)
$ws.Range("G2").Value = 1

# Row 3 - Austin: increment May 2021 (G) and May 2017 (AA)
$ws.Range("G3").Value = 3
$ws.Range("AA3").Value = 3

# Row 4 - Humboldt Park: new May 2022 (B) value, increment May 2021 (G)
$ws.Range("B4").Value = 1
$ws.Range("G4").Value = 3

# Row 5 - Garfield Park: increment April 2022 (C)
$ws.Range("C5").Value = 7

# Row 8 - South Shore: new May 2022 (B) value, increment May 2015 (AK)
$ws.Range("B8").Value = 1
$ws.Range("AK8").Value = 2

# Row 13 - Washington Heights: new May 2017 (AA) and May 2015 (AK) values
$ws.Range("AA13").Value = 1
$ws.Range("AK13").Value = 1

# Row 20 - Woodlawn: new May 2021 (G) value
$ws.Range("G20").Value = 1

# Row 21 - Chatham: new May 2022 (B) value
$ws.Range("B21").Value = 1

# Row 25 - Auburn Gresham: new May 2022 (B) value
$ws.Range("B25").Value = 1

# Row 30 - West Loop: new May 2022 (B) value
$ws.Range("B30").Value = 1

# Row 38 - Douglas: increment May 2021 (G)
$ws.Range("G38").Value = 2
